# "add a couple of more pages"
# The underlying dataset dropped two records (original rows 2 and 6, i.e.
# Id 88353 and Id 88357). Every row below each deleted record shifts up,
# the sheet's used range shrinks from A1:AX11 to A1:AX9, a handful of
# columns pick up explicit widths, and the sheet view's selection moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two dropped records. Row 2 (Id 88353) goes first; once it is
# gone, the old row 6 (Id 88357) has shifted up to row 5, so deleting row 5
# next removes exactly that record. Everything else cascades upward
# automatically, which reproduces the new Id/metric alignment seen in the
# diff without having to re-type every shifted cell by hand.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(5).Delete()

# New explicit column widths introduced in this revision.
$ws.Columns.Item(4).ColumnWidth = 28.166666666666668   # D  -> stored width 29
$ws.Columns.Item(5).ColumnWidth = 39.0                  # E  -> stored width 39.83203125
$ws.Columns.Item(6).ColumnWidth = 15.333333333333334    # F  -> stored width 16.1640625
$ws.Columns.Item(7).ColumnWidth = 13.833333333333334    # G  -> stored width 14.6640625
$ws.Columns.Item(12).ColumnWidth = 17.833333333333332   # L  -> stored width 18.6640625
$ws.Columns.Item(13).ColumnWidth = 20.333333333333332   # M  -> stored width 21.1640625
$ws.Columns.Item(14).ColumnWidth = 16.166666666666668   # N  -> stored width 17
$ws.Columns.Item(15).ColumnWidth = 44.833333333333336   # O  -> stored width 45.6640625
$ws.Columns.Item(16).ColumnWidth = 15.833333333333334   # P  -> stored width 16.6640625
$ws.Columns.Item(21).ColumnWidth = 21.0                 # U  -> stored width 21.83203125

# Scroll the view over to column O and move the active selection, matching
# the refreshed sheetView in the saved workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 15
$win.ScrollRow = 1
$ws.Range("E15").Select()
